$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.381.87"
$ws.Range("E2").Value = "  -1.37%  "
$ws.Range("D3").Value = "'1.711.00"
$ws.Range("E3").Value = "  -1.55%  "
$ws.Range("D4").Value = "'1.005"
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").Value = "'224.43"
$ws.Range("E5").Value = "  -1.43%  "
$ws.Range("D6").Value = "'0.5334"
$ws.Range("E6").Value = "  -2.31%  "
$ws.Range("D7").Value = "'1.005"
$ws.Range("E7").Value = "  +0.19%  "
$ws.Range("D8").Value = "'0.2668"
$ws.Range("E8").Value = "  -3.37%  "
$ws.Range("D9").Value = "'0.06611"
$ws.Range("E9").Value = "  -1.52%  "
$ws.Range("D10").Value = "'20.92"
$ws.Range("E10").Value = "  -4.55%  "
$ws.Range("D11").Value = "'0.07643"
$ws.Range("E11").Value = "  -1.76%  "
$ws.Range("D12").Value = "'4.557"
$ws.Range("E12").Value = "  -2.89%  "
$ws.Range("D13").Value = "'1.733.83"
$ws.Range("E13").Value = "  -1.12%  "
$ws.Range("D14").Value = "'1.948.22"
$ws.Range("E14").Value = "  -1.42%  "
$ws.Range("D15").Value = "'0.5765"
$ws.Range("E15").Value = "  -3.78%  "
$ws.Range("D16").Value = "'0.0₅8174"
$ws.Range("E16").Value = "  -3.07%  "
$ws.Range("D17").Value = "'67.85"
$ws.Range("E17").Value = "  -2.33%  "
$ws.Range("D18").Value = "'27.405.93"
$ws.Range("E18").Value = "  -1.30%  "
$ws.Range("D19").Value = "'216.25"
$ws.Range("E19").Value = "  -4.68%  "
$ws.Range("E20").Value = "  +0.15%  "
$ws.Range("D21").Value = "'4.663"
$ws.Range("E21").Value = "  -3.68%  "
$ws.Range("E22").Value = "  -4.07%  "
$ws.Range("D23").Value = "'5.983"
$ws.Range("E23").Value = "  -4.00%  "
$ws.Range("D24").Value = "'1.006"
$ws.Range("E24").Value = "  +0.17%  "
$ws.Range("D25").Value = "'142.49"
$ws.Range("E25").Value = "  -3.16%  "
$ws.Range("D26").Value = "'1.733"
$ws.Range("E26").Value = "  +1.50%  "
$ws.Range("D27").Value = "'0.1216"
$ws.Range("E27").Value = "  -2.91%  "
$ws.Range("E28").Value = "  -2.39%  "
$ws.Range("D29").Value = "'16.26"
$ws.Range("E29").Value = "  -5.34%  "
$ws.Range("D30").Value = "'0.05403"
$ws.Range("E30").Value = "  -4.83%  "
$ws.Range("D31").Value = "'1.292"
$ws.Range("E31").Value = "  -1.69%  "
$ws.Range("D32").Value = "'3.488"
$ws.Range("E32").Value = "  -5.78%  "
$ws.Range("D33").Value = "'3.427"
$ws.Range("E33").Value = "  -2.60%  "
$ws.Range("D34").Value = "'1.643"
$ws.Range("E34").Value = "  -2.60%  "
$ws.Range("D35").Value = "'2.882"
$ws.Range("E35").Value = "  +0.96%  "
$ws.Range("D36").Value = "'0.9490"
$ws.Range("E36").Value = "  -2.80%  "
$ws.Range("D37").Value = "'2.418"
$ws.Range("E37").Value = "  -1.28%  "
$ws.Range("D38").Value = "'0.5836"
$ws.Range("E38").Value = "  -2.23%  "
$ws.Range("D39").Value = "'0.01633"
$ws.Range("E39").Value = "  -2.23%  "
$ws.Range("D40").Value = "'5.861"
$ws.Range("E40").Value = "  -0.77%  "
$ws.Range("D41").Value = "'1.046.28"
$ws.Range("E41").Value = "  -0.52%  "
$ws.Range("D42").Value = "'1.006"
$ws.Range("E42").Value = "  +0.20%  "
$ws.Range("D43").Value = "'0.8416"
$ws.Range("E43").Value = "  -0.78%  "
$ws.Range("D44").Value = "'100.77"
$ws.Range("E44").Value = "  -1.38%  "
$ws.Range("D45").Value = "'1.855.65"
$ws.Range("E45").Value = "  -1.33%  "
$ws.Range("E46").Value = "  +2.33%  "
$ws.Range("D47").Value = "'58.00"
$ws.Range("E47").Value = "  -2.59%  "
$ws.Range("D48").Value = "'0.4521"
$ws.Range("E48").Value = "  +1.86%  "
$ws.Range("D49").Value = "'1.004"
$ws.Range("E49").Value = "  +0.15%  "
$ws.Range("D50").Value = "'8.056"
$ws.Range("E50").Value = "  -2.70%  "
$ws.Range("D51").Value = "'0.05242"
$ws.Range("E51").Value = "  -1.38%  "
